$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 22, pushing the existing row 22 (and below) down to row 23.
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new weekly price record.
$ws.Range("A22").Value = 10
$ws.Range("B22").Value = "Vega Modelo de Temuco"
$ws.Range("C22").Value = "La Araucanía"
$ws.Range("D22").Value = 44448
$ws.Range("D22").NumberFormat = $ws.Range("D23").NumberFormat
$ws.Range("E22").Value = 9
$ws.Range("F22").Value = 100112026
$ws.Range("G22").Value = "Haba"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 35
$ws.Range("K22").Value = 15000
$ws.Range("L22").Value = 15000
$ws.Range("M22").Value = 15000
$ws.Range("N22").Value = "$/saco 25 kilos"
$ws.Range("O22").Value = "Provincia de Limarí"
$ws.Range("P22").Value = 600
$ws.Range("Q22").Value = 25
$ws.Range("R22").Value = "Hortaliza"
